# Refresh the cryptocurrency Price (column D) and Volume(1h) (column E)
# figures in the tracker sheet with the latest scraped snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '63.375.34' },
    @{ Cell = 'E2'; Value = '  -3.75%  ' },
    @{ Cell = 'D3'; Value = '2.585.41' },
    @{ Cell = 'E3'; Value = '  -2.98%  ' },
    @{ Cell = 'E4'; Value = '  +0.05%  ' },
    @{ Cell = 'D5'; Value = '571.68' },
    @{ Cell = 'E5'; Value = '  -4.50%  ' },
    @{ Cell = 'D6'; Value = '154.55' },
    @{ Cell = 'E6'; Value = '  -2.43%  ' },
    @{ Cell = 'E7'; Value = '  +0.10%  ' },
    @{ Cell = 'E8'; Value = '  -4.78%  ' },
    @{ Cell = 'E9'; Value = '  -7.47%  ' },
    @{ Cell = 'D10'; Value = '5.83' },
    @{ Cell = 'E10'; Value = '  -0.44%  ' },
    @{ Cell = 'D11'; Value = '0.379' },
    @{ Cell = 'E11'; Value = '  -5.90%  ' },
    @{ Cell = 'E12'; Value = '  -0.38%  ' },
    @{ Cell = 'D13'; Value = '28.14' },
    @{ Cell = 'E13'; Value = '  -3.02%  ' },
    @{ Cell = 'D14'; Value = '3.056.95' },
    @{ Cell = 'E14'; Value = '  -2.78%  ' },
    @{ Cell = 'D16'; Value = '63.245.18' },
    @{ Cell = 'E16'; Value = '  -3.72%  ' },
    @{ Cell = 'D17'; Value = '2.616.01' },
    @{ Cell = 'E17'; Value = '  -2.62%  ' },
    @{ Cell = 'D18'; Value = '11.93' },
    @{ Cell = 'E18'; Value = '  -5.36%  ' },
    @{ Cell = 'D19'; Value = '7.48' },
    @{ Cell = 'E19'; Value = '  -0.43%  ' },
    @{ Cell = 'D20'; Value = '4.53' },
    @{ Cell = 'E20'; Value = '  -5.90%  ' },
    @{ Cell = 'D21'; Value = '341.36' },
    @{ Cell = 'E21'; Value = '  -2.92%  ' },
    @{ Cell = 'E22'; Value = '  +0.07%  ' },
    @{ Cell = 'D23'; Value = '67.31' },
    @{ Cell = 'E23'; Value = '  -3.64%  ' },
    @{ Cell = 'E24'; Value = '  -0.42%  ' },
    @{ Cell = 'E25'; Value = '  -4.12%  ' },
    @{ Cell = 'D26'; Value = '9.09' },
    @{ Cell = 'E26'; Value = '  -5.87%  ' },
    @{ Cell = 'D27'; Value = '575.87' },
    @{ Cell = 'E27'; Value = '  +0.89%  ' },
    @{ Cell = 'E28'; Value = '  -3.99%  ' },
    @{ Cell = 'D29'; Value = '0.999' },
    @{ Cell = 'E29'; Value = '  -0.02%  ' },
    @{ Cell = 'E30'; Value = '  -1.69%  ' },
    @{ Cell = 'D31'; Value = '7.86' },
    @{ Cell = 'E31'; Value = '  -4.13%  ' },
    @{ Cell = 'E32'; Value = '  -5.05%  ' },
    @{ Cell = 'D33'; Value = '1.72' },
    @{ Cell = 'E33'; Value = '  -5.65%  ' },
    @{ Cell = 'D34'; Value = '6.50' },
    @{ Cell = 'E34'; Value = '  -3.21%  ' },
    @{ Cell = 'D35'; Value = '5.43' },
    @{ Cell = 'E35'; Value = '  -2.62%  ' },
    @{ Cell = 'D37'; Value = '0.998' },
    @{ Cell = 'E37'; Value = '  -0.07%  ' },
    @{ Cell = 'D38'; Value = '19.69' },
    @{ Cell = 'E38'; Value = '  -4.69%  ' },
    @{ Cell = 'D39'; Value = '154.19' },
    @{ Cell = 'E40'; Value = '  -5.27%  ' },
    @{ Cell = 'E41'; Value = '  -0.02%  ' },
    @{ Cell = 'D42'; Value = '41.25' },
    @{ Cell = 'D43'; Value = '2.45' },
    @{ Cell = 'E43'; Value = '  +5.68%  ' },
    @{ Cell = 'D44'; Value = '155.76' },
    @{ Cell = 'E44'; Value = '  -3.83%  ' },
    @{ Cell = 'D45'; Value = '3.90' },
    @{ Cell = 'E45'; Value = '  -5.21%  ' },
    @{ Cell = 'D46'; Value = '23.02' },
    @{ Cell = 'E46'; Value = '  -0.49%  ' },
    @{ Cell = 'D47'; Value = '0.0586' },
    @{ Cell = 'E47'; Value = '  -5.41%  ' },
    @{ Cell = 'D48'; Value = '0.624' },
    @{ Cell = 'E48'; Value = '  -3.19%  ' },
    @{ Cell = 'E49'; Value = '  -2.30%  ' },
    @{ Cell = 'D50'; Value = '0.0245' },
    @{ Cell = 'E50'; Value = '  -5.17%  ' },
    @{ Cell = 'E51'; Value = '  -5.52%  ' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $row = $cell.Row

    # Force text entry so numeric-looking strings (e.g. "571.68", "0.0586")
    # are not silently reinterpreted as numbers by Excel, matching the
    # plain-text values already used throughout the sheet.
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value

    # Re-align the cell style with the rest of the (unstyled) data rows so
    # forcing text formatting above doesn't leave a stray quote-prefix /
    # text-format style behind; column B never carries an explicit style.
    $cell.Style = $ws.Cells.Item($row, 2).Style
}
